# Update 2p0. Convention change to support multi-axle vehicles
#
# - Adds a new "Truck_Amandla" sheet (cloned from Trailer_Elula) positioned
#   right after Bus_Makhulu.
# - Adds a new "Trailer_Kumanzi" sheet (cloned from Trailer_Thwala)
#   positioned right after Trailer_Thwala (i.e. at the end).
# - Updates the per-sheet "Instance" (H3), CD (H6) and force coefficients
#   (F9:H9) on the two new sheets.

$wb = $excel.ActiveWorkbook

# --- Truck_Amandla: clone of Trailer_Elula, placed after Bus_Makhulu -------
$busMakhulu  = $wb.Worksheets.Item("Bus_Makhulu")
$trailerElula = $wb.Worksheets.Item("Trailer_Elula")
$trailerElula.Copy($null, $busMakhulu)

$truck = $wb.Worksheets.Item("Trailer_Elula (2)")
$truck.Name = "Truck_Amandla"

$truck.Range("H3").Value = "Truck_Amandla"
$truck.Range("H6").Value = 0.43
$truck.Range("F9").Value = -1.2
$truck.Range("G9").Value = 0
$truck.Range("H9").Value = 1.1

$truck.Range("H5:H9").Select()

# --- Trailer_Kumanzi: clone of Trailer_Thwala, placed at the end -----------
$trailerThwala = $wb.Worksheets.Item("Trailer_Thwala")
$trailerThwala.Copy($null, $trailerThwala)

$kumanzi = $wb.Worksheets.Item("Trailer_Thwala (2)")
$kumanzi.Name = "Trailer_Kumanzi"

$kumanzi.Range("H3").Value = "Trailer_Kumanzi"
$kumanzi.Range("H6").Value = 0.43
$kumanzi.Range("F9").Value = 5
$kumanzi.Range("G9").Value = 0
$kumanzi.Range("H9").Value = 2

$kumanzi.Range("J20").Select()
$kumanzi.Activate()
